# Add new columns I ("I0") and J ("IF") to the sheet, mirroring the
# style/structure of the existing header + data columns (A..H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1), using the same style as the other header cells (B1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2..65, column I ("I0") and column J ("IF").
$iValues = @(6,7,6,5,8,8,5,7,9,9,5,8,8,8,9,9,7,8,6,8,11,7,9,6,7,8,10,7,7,8,9,7,4,9,7,7,6,8,6,8,7,7,8,6,6,9,8,7,9,10,6,6,7,7,8,9,2,5,4,2,5,4,4,2)
$jValues = @(7,7,6,5,8,8,6,7,9,9,5,8,8,8,9,9,7,8,6,8,11,8,9,7,7,8,10,7,7,8,9,7,4,9,7,8,6,9,6,8,7,7,8,6,6,9,8,7,9,10,6,6,7,7,8,9,2,5,5,2,5,4,4,2)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
